# Add extra columns (WIN, TOP4, TOP5, RELEGATION) before ExpPoints, re-sort the
# table by the refreshed ExpPoints values and move ExpPoints into column G,
# leaving the new C:F columns blank placeholders for the upcoming Monte Carlo
# simulation percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data, already sorted descending by the refreshed ExpPoints values.
$data = @(
    @("Arsenal", 78.16550934895004),
    @("Liverpool", 72.90942580624936),
    @("Manchester City", 69.46422848523719),
    @("Chelsea", 63.12886895979322),
    @("Crystal Palace", 58.78560160382771),
    @("Aston Villa", 55.51762574508723),
    @("Newcastle United", 55.36782147208676),
    @("Brighton & Hove Albion", 54.35772172953664),
    @("AFC Bournemouth", 54.08348221059467),
    @("Tottenham Hotspur", 53.83456459564835),
    @("Manchester United", 51.64669391800064),
    @("Brentford", 50.0761895883378),
    @("Everton", 45.67963808229034),
    @("Fulham", 42.62410787626668),
    @("Nottingham Forest", 41.12453832570574),
    @("Sunderland", 39.76161622693008),
    @("West Ham United", 38.22778832823779),
    @("Leeds United", 36.26706835306463),
    @("Burnley", 34.79629067768817),
    @("Wolverhampton Wanderers", 31.47489508579568)
)

# Move the ExpPoints header from C1 to G1, and insert the new headers in
# between.
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP4"
$ws.Range("E1").Value = "TOP5"
$ws.Range("F1").Value = "RELEGATION"
$ws.Range("G1").Value = "ExpPoints"

# Match the bold / centered / boxed header formatting already used by A1:B1.
$headerRng = $ws.Range("C1:G1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160
$headerRng.Borders.LineStyle = 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $team = $data[$i][0]
    $points = $data[$i][1]

    $ws.Cells.Item($row, 2).Value = $team   # B: Team (re-sorted)

    # C:F are blank placeholder cells (WIN / TOP4 / TOP5 / RELEGATION) to be
    # filled in later by the Monte Carlo simulation. Touching NumberFormat
    # after blanking keeps a real (empty) cell on the sheet instead of
    # letting it get swept away as unused.
    foreach ($col in 3..6) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = ""
        $cell.NumberFormat = "General"
    }

    $ws.Cells.Item($row, 7).Value = $points # G: ExpPoints (refreshed value)
}
